$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column D (Price) and E (Volume(1h)) values for rows with changed data ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.817.11"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.345.20"
$ws.Range("E3").Value = "  -1.13%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.33"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.670"
$ws.Range("E6").Value = "  -3.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.34"
$ws.Range("E7").Value = "  -5.64%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.592"
$ws.Range("E9").Value = "  -6.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0997"
$ws.Range("E10").Value = "  -3.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.44"
$ws.Range("E11").Value = "  +1.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "32.49"
$ws.Range("E12").Value = "  -2.08%  "
$ws.Range("E13").Value = "  -0.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.12"
$ws.Range("E14").Value = "  -6.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.694.25"
$ws.Range("E15").Value = "  -0.93%  "
$ws.Range("E16").Value = "  -4.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.901"
$ws.Range("E17").Value = "  -3.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.346.10"
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.706.31"
$ws.Range("E19").Value = "  -4.78%  "
$ws.Range("E20").Value = "  -1.15%  "
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.77"
$ws.Range("E23").Value = "  -1.74%  "
$ws.Range("E24").Value = "  +8.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("E26").Value = "  +1.21%  "
$ws.Range("E27").Value = "  -2.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.39"
$ws.Range("E28").Value = "  -7.65%  "
$ws.Range("E29").Value = "  -1.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "175.96"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.29"
$ws.Range("E31").Value = "  -4.43%  "
$ws.Range("E32").Value = "  -1.55%  "
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0739"
$ws.Range("E34").Value = "  -2.00%  "
$ws.Range("E35").Value = "  -4.89%  "
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.74"
$ws.Range("E37").Value = "  -2.96%  "
$ws.Range("E38").Value = "  -3.95%  "
$ws.Range("E39").Value = "  -1.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0272"
$ws.Range("E40").Value = "  -1.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "65.41"
$ws.Range("E41").Value = "  +19.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.23"
$ws.Range("E42").Value = "  +15.57%  "
$ws.Range("E43").Value = "  +3.20%  "
$ws.Range("E44").Value = "  +6.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.90"
$ws.Range("E45").Value = "  -1.07%  "
$ws.Range("E46").Value = "  -2.19%  "
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("E48").Value = "  -2.72%  "
$ws.Range("E49").Value = "  -2.99%  "

# --- Rows 50 and 51: content swap (ARBITRUM moves up, Aave moves down) with new values ---
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.15"
$ws.Range("E50").Value = "  -4.66%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "98.37"
$ws.Range("E51").Value = "  -4.61%  "
